$p = $ppt.ActivePresentation

# Slide 5 (sldId 268): update title text on shape id 86
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(1)
$sh5.TextFrame.TextRange.Text = "[Component 1] (Instruction&yes/no checker)"

# Slide 11 (sldId 270): move title shape (id 2, "Title 1") up
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(1)
$sh11.Top = 300470 / 12700
